$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain plain text
# (matching the source columns inlineStr/text formatting), so force
# a Text number format on them before assigning the value.
$textCells = @("D5", "D6", "D12", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D31", "D32", "D36", "D38", "D40", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range("D2").Value = "95.796.94"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "3.552.30"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "238.71"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "650.86"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  +10.25%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +6.41%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "3.550.96"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "43.00"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "6.36"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "4.211.50"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "95.686.26"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "3.556.52"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "7.75"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("D20").Value = "12.41"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "17.63"
$ws.Range("E21").Value = "  -1.38%  "
$ws.Range("D22").Value = "0.517"
$ws.Range("E22").Value = "  +6.81%  "
$ws.Range("D23").Value = "502.04"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "3.37"
$ws.Range("E24").Value = "  -6.63%  "
$ws.Range("D25").Value = "6.85"
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "95.41"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").Value = "12.66"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").Value = "3.743.22"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  +9.73%  "
$ws.Range("D31").Value = "2.98"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("D32").Value = "11.28"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").Value = "31.18"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("E37").Value = "  +6.86%  "
$ws.Range("D38").Value = "607.32"
$ws.Range("E38").Value = "  +6.41%  "
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "1.60"
$ws.Range("E40").Value = "  +8.49%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("E44").Value = "  +5.28%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "23.50"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0417"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "33.73"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "8.10"
$ws.Range("E51").Value = "  +0.85%  "
